$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.739.42"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").Value = "1.873.94"
$ws.Range("E3").Value = "  -2.08%  "

$ws.Range("E4").Value = "  -0.95%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.23"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.685"
$ws.Range("E6").Value = "  -2.80%  "

$ws.Range("E7").Value = "  -1.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.51"
$ws.Range("E8").Value = "  +1.99%  "

$ws.Range("E9").Value = "  -3.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "50.69"
$ws.Range("E10").Value = "  -3.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0737"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0967"
$ws.Range("E12").Value = "  -2.52%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.85"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.143.50"
$ws.Range("E14").Value = "  -2.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.713"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.88"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "1.872.74"
$ws.Range("E17").Value = "  -2.17%  "

$ws.Range("D18").Value = "34.702.99"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.77"
$ws.Range("E19").Value = "  -0.69%  "

$ws.Range("D20").Value = "0.0₃0821"
$ws.Range("E20").Value = "  -0.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "246.68"
$ws.Range("E21").Value = "  +1.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.72"
$ws.Range("E22").Value = "  -3.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.91"
$ws.Range("E23").Value = "  -3.19%  "

$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +2.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -4.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.91"
$ws.Range("E27").Value = "  -2.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.38"
$ws.Range("E28").Value = "  -3.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.22"
$ws.Range("E29").Value = "  -2.97%  "

$ws.Range("E30").Value = "  -4.02%  "

$ws.Range("D31").Value = "4.128.62"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("B32").Value = "TrustWalletToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.64"
$ws.Range("E32").Value = "  +11.19%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.25"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0579"
$ws.Range("E34").Value = "  +0.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.16"
$ws.Range("E35").Value = "  -1.35%  "

$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  -5.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.832"
$ws.Range("E38").Value = "  -8.70%  "

$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.15"
$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.20"
$ws.Range("E41").Value = "  -1.27%  "

$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("E44").Value = "  -5.30%  "

$ws.Range("D45").Value = "1.291.09"
$ws.Range("E45").Value = "  -4.55%  "

$ws.Range("E46").Value = "  -4.48%  "

$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.06"
$ws.Range("E49").Value = "  -1.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0761"
$ws.Range("E50").Value = "  +5.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.46"
$ws.Range("E51").Value = "  -1.44%  "
